$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) data for each named row (rows 2-13), keyed by
# the method name in column A.
$data = @{}
for ($r = 2; $r -le 13; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    $data[$name] = @(
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 5).Value2,
        $ws.Cells.Item($r, 6).Value2
    )
}

# Fix up the metrics for KAOGExp — computing metrics for results without
# optimization (handling None values that previously produced incomplete
# metrics).
$data["KAOGExp"] = @(1, 17.75, 12.93879281845967, 12.7541939361204, 1)

# New row order (rows 2-13) after recomputing metrics / reshuffling results.
$newOrder = @("KAOGExp", "cruds", "wachter", "face-knn", "cem", "cem-vae", "dice", "face-epsilon", "clue", "ar", "cchvae", "gs")

$row = 2
foreach ($name in $newOrder) {
    $ws.Cells.Item($row, 1).Value = $name
    $vals = $data[$name]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $row++
}
